$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 15627463
$ws.Range("I33").Value = 20834034
$ws.Range("K33").Value = 20834034
$ws.Range("M33").Value = -20833805
$ws.Range("H62").Value = 823.75
$ws.Range("I62").Value = 823.75
$ws.Range("K62").Value = 823.75
$ws.Range("M62").Value = -199.75
$ws.Range("H64").Value = 83339680
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H65").Value = 823.75
$ws.Range("I65").Value = 823.75
$ws.Range("K65").Value = 4118.75
$ws.Range("M65").Value = -998.75
$ws.Range("H67").Value = 83339680
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H76").Value = 7136.6665
$ws.Range("J76").Value = 10560.5
$ws.Range("L76").Value = 10560.5
$ws.Range("N76").Value = -11190.5
$ws.Range("H79").Value = 7136.6665
$ws.Range("J79").Value = 10560.5
$ws.Range("L79").Value = 10560.5
$ws.Range("N79").Value = -12744.5
$ws.Range("H86").Value = 7237.7856
$ws.Range("I86").Value = 6433.4
$ws.Range("J86").Value = 9248.75
$ws.Range("K86").Value = 6433.4
$ws.Range("L86").Value = 9248.75
$ws.Range("M86").Value = -5310.4
$ws.Range("N86").Value = -11494.75
$ws.Range("H89").Value = 7237.7856
$ws.Range("I89").Value = 6433.4
$ws.Range("J89").Value = 9248.75
$ws.Range("K89").Value = 32167
$ws.Range("L89").Value = 46243.75
$ws.Range("M89").Value = -26551
$ws.Range("N89").Value = -57475.75
$ws.Range("H116").Value = 4484.857
$ws.Range("I116").Value = 2898.5
$ws.Range("K116").Value = 2898.5
$ws.Range("M116").Value = 543.5
$ws.Range("H137").Value = 10812.12
$ws.Range("I137").Value = 3599.5454
$ws.Range("J137").Value = 16479.143
$ws.Range("K137").Value = 10798.6362
$ws.Range("L137").Value = 49437.429
$ws.Range("M137").Value = -8248.636200000001
$ws.Range("N137").Value = -54537.429
$ws.Range("H141").Value = 4179.625
$ws.Range("I141").Value = 4426.4614
$ws.Range("K141").Value = 13279.3842
$ws.Range("M141").Value = -8099.3842

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 33666.668
$ws.Range("I43").Value = 32000
$ws.Range("K43").Value = 32000
$ws.Range("M43").Value = -31687
$ws.Range("H61").Value = 95489.24000000001
$ws.Range("I61").Value = 2716.1667
$ws.Range("K61").Value = 2716.1667
$ws.Range("M61").Value = -2504.1667
$ws.Range("H62").Value = 88249
$ws.Range("J62").Value = 88249
$ws.Range("L62").Value = 88249
$ws.Range("N62").Value = -89497
$ws.Range("H65").Value = 88249
$ws.Range("J65").Value = 88249
$ws.Range("L65").Value = 264747
$ws.Range("N65").Value = -270987
$ws.Range("H132").Value = 8364708
$ws.Range("I132").Value = 18404
$ws.Range("K132").Value = 55212
$ws.Range("M132").Value = -52682
$ws.Range("H136").Value = 95489.24000000001
$ws.Range("I136").Value = 2716.1667
$ws.Range("K136").Value = 8148.500100000001
$ws.Range("M136").Value = -5598.500100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 23816.186
$ws.Range("I99").Value = 22480.39
$ws.Range("K99").Value = 22480.39
$ws.Range("M99").Value = -20982.39
$ws.Range("H134").Value = 42864.805
$ws.Range("I134").Value = 41554.89
$ws.Range("K134").Value = 124664.67
$ws.Range("M134").Value = -122129.67

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22220.27
$ws.Range("I31").Value = 12199.8
$ws.Range("K31").Value = 12199.8
$ws.Range("M31").Value = -11904.8
$ws.Range("H34").Value = 22220.27
$ws.Range("I34").Value = 12199.8
$ws.Range("K34").Value = 12199.8
$ws.Range("M34").Value = -11997.8
$ws.Range("H58").Value = 13235.475
$ws.Range("I58").Value = 4909.129
$ws.Range("K58").Value = 4909.129
$ws.Range("M58").Value = -4706.129
$ws.Range("H63").Value = 14519.091
$ws.Range("I63").Value = 7000
$ws.Range("K63").Value = 7000
$ws.Range("H66").Value = 14519.091
$ws.Range("I66").Value = 7000
$ws.Range("K66").Value = 21000
$ws.Range("H86").Value = 11245.389
$ws.Range("I86").Value = 12043.6
$ws.Range("J86").Value = 7254.3335
$ws.Range("K86").Value = 12043.6
$ws.Range("L86").Value = 7254.3335
$ws.Range("M86").Value = -10920.6
$ws.Range("N86").Value = -9500.333500000001
$ws.Range("H89").Value = 11245.389
$ws.Range("I89").Value = 12043.6
$ws.Range("J89").Value = 7254.3335
$ws.Range("K89").Value = 60218
$ws.Range("L89").Value = 36271.6675
$ws.Range("M89").Value = -54602
$ws.Range("N89").Value = -47503.6675
$ws.Range("H134").Value = 24395306
$ws.Range("I134").Value = 1752.6897
$ws.Range("J134").Value = 83346390
$ws.Range("K134").Value = 5258.0691
$ws.Range("L134").Value = 250039170
$ws.Range("M134").Value = -2723.0691
$ws.Range("N134").Value = -250044240
$ws.Range("H136").Value = 13235.475
$ws.Range("I136").Value = 4909.129
$ws.Range("K136").Value = 14727.387
$ws.Range("M136").Value = -12177.387
$ws.Range("M63").Value = -6314
$ws.Range("M66").Value = -17568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 66668500
$ws.Range("J32").Value = 33333666
$ws.Range("L32").Value = 100000998
$ws.Range("N32").Value = -100001564
$ws.Range("H97").Value = 322.66666
$ws.Range("J97").Value = 750
$ws.Range("L97").Value = 2250
$ws.Range("N97").Value = -3242

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3791676.5
$ws.Range("I102").Value = 6584154.5
$ws.Range("K102").Value = 6584154.5
$ws.Range("M102").Value = -6582532.5
$ws.Range("H132").Value = 15247.2
$ws.Range("I132").Value = 13672.333
$ws.Range("J132").Value = 15922.143
$ws.Range("K132").Value = 41016.999
$ws.Range("L132").Value = 47766.429
$ws.Range("M132").Value = -38486.999
$ws.Range("N132").Value = -52826.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 44747.5
$ws.Range("I45").Value = 42996.668
$ws.Range("K45").Value = 42996.668
$ws.Range("M45").Value = -42589.668
$ws.Range("H61").Value = 2679.4285
$ws.Range("I61").Value = 2053.25
$ws.Range("J61").Value = 4683.2
$ws.Range("K61").Value = 2053.25
$ws.Range("L61").Value = 4683.2
$ws.Range("M61").Value = -1851.25
$ws.Range("N61").Value = -5087.2
$ws.Range("H113").Value = 2679.4285
$ws.Range("I113").Value = 2053.25
$ws.Range("J113").Value = 4683.2
$ws.Range("K113").Value = 2053.25
$ws.Range("L113").Value = 4683.2
$ws.Range("M113").Value = 116.75
$ws.Range("N113").Value = -9023.200000000001
$ws.Range("H122").Value = 20995422
$ws.Range("I122").Value = 40062020
$ws.Range("J122").Value = 1928824.6
$ws.Range("K122").Value = 120186060
$ws.Range("L122").Value = 5786473.800000001
$ws.Range("M122").Value = -120183610
$ws.Range("N122").Value = -5791373.800000001
$ws.Range("H132").Value = 3094865.5
$ws.Range("I132").Value = 1916
$ws.Range("K132").Value = 5748
$ws.Range("M132").Value = -3218

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H100").Value = 753
$ws.Range("I100").Value = 690
$ws.Range("J100").Value = 784.5
$ws.Range("K100").Value = 1380
$ws.Range("L100").Value = 1569
$ws.Range("M100").Value = -839
$ws.Range("N100").Value = -2651
$ws.Range("N98").ClearContents()
